# Update the "Report Generated On" timestamp in D5
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:49 PM"

# Zero out the Total Billed Amount and all daily pricing / total cells in column H
$cellsToZero = @("C8", "H16", "H17", "H18", "H19", "H20", "H25", "H26", "H27", "H28", "H29", "H30", "H31", "H36", "H37", "H38", "H39", "H40", "H41", "H42", "H47", "H48", "H49", "H50", "H51")

foreach ($addr in $cellsToZero) {
    $ws.Range($addr).Value = 0
}
